# Adds new planet-selection/customization requirements (SYS-PLT-001..006)
# to the MASTER SPREADSHEET sheet, rows 271-277, and updates the sheet's
# saved view (scroll position / zoom / selection) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MASTER SPREADSHEET")
$ws.Activate()

# Reference cells whose existing style encodes the urgency color-coding
# used throughout the sheet (HIGH = red fill, MID = yellow fill,
# LOW = green fill). We copy their formatting onto the new urgency
# cells instead of inventing new style indices.
$highRef = $ws.Range("C265")
$midRef  = $ws.Range("C269")
$lowRef  = $ws.Range("C212")

function Set-Urgency($addr, $refCell, $text) {
    $refCell.Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range($addr).Value = $text
}

$rows = @(
    @{ Row = 271; Id = "SYS-PLT-001";       Urgency = "HIGH"; Ref = $highRef; Desc = "The program shall have a defined template for planets." },
    @{ Row = 272; Id = "SYS-PLT-002";       Urgency = "MID";  Ref = $midRef;  Desc = "The program shall allow for planets to be referenced externally" },
    @{ Row = 273; Id = "SYS-PLT-003-001";   Urgency = "HIGH"; Ref = $highRef; Desc = "The program shall have a 'property' and 'function' element included in the planet definition." },
    @{ Row = 274; Id = "SYS-PLT-003-002";   Urgency = "HIGH"; Ref = $highRef; Desc = "The program shall present those elements only once." },
    @{ Row = 275; Id = "SYS-PLT-004";       Urgency = "HIGH"; Ref = $highRef; Desc = "The program shall have the planet defintion that includes the following: semimajor axis, semiminor axis, equalatioral axis, equalatorial radius, polar radius, rotation rate, gravity, and J2." },
    @{ Row = 276; Id = "SYS-PLT-005";       Urgency = "LOW";  Ref = $lowRef;  Desc = "The program shall allow for planets to be named by the user." },
    @{ Row = 277; Id = "SYS-PLT-006";       Urgency = "HIGH"; Ref = $highRef; Desc = "The program must have the J2 element as a POSITIVE number." }
)

foreach ($r in $rows) {
    $ws.Range("B$($r.Row)").Value = $r.Id
    Set-Urgency "C$($r.Row)" $r.Ref $r.Urgency
    $ws.Range("D$($r.Row)").Value = $r.Desc
}

# Update the saved view to match where the author ended up after adding
# the rows above.
$ws.Range("D274").Select() | Out-Null
$win = $excel.ActiveWindow
$win.Zoom = 85
$win.ScrollRow = 264
$win.ScrollColumn = 1
